$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("itself to align", $true, $false, $false, $false, $false, $true, 1, $false, "itself (e.g. x/y2=x/y1*sin?/cos?(a)) to align", 2)
